# Apply scheduled-runner profit/price updates to the Leve profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 8502.5  # H7 was 10005
$ws.Cells.Item(7, 10).Value = 7000  # J7 was 0
$ws.Cells.Item(7, 12).Value = 7000  # L7 was 0
$ws.Cells.Item(7, 14).Value = -7224  # N7 was empty

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(14, 8).Value = 8502.5  # H14 was 10005
$ws.Cells.Item(14, 10).Value = 7000  # J14 was 0
$ws.Cells.Item(14, 12).Value = 7000  # L14 was 0
$ws.Cells.Item(14, 14).Value = -7382  # N14 was empty

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 617.1429000000001  # H107 was 393.33334
$ws.Cells.Item(107, 9).Value = 426  # I107 was 393.33334
$ws.Cells.Item(107, 10).Value = 1095  # J107 was 0
$ws.Cells.Item(107, 11).Value = 426  # K107 was 393.33334
$ws.Cells.Item(107, 12).Value = 1095  # L107 was 0
$ws.Cells.Item(107, 13).Value = 1494  # M107 was 1526.66666
$ws.Cells.Item(107, 14).Value = -4935  # N107 was empty

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 1986.0667  # H116 was 1852.6923
$ws.Cells.Item(116, 9).Value = 1798  # I116 was 1609.4445
$ws.Cells.Item(116, 10).Value = 2080.1  # J116 was 2400
$ws.Cells.Item(116, 11).Value = 1798  # K116 was 1609.4445
$ws.Cells.Item(116, 12).Value = 2080.1  # L116 was 2400
$ws.Cells.Item(116, 13).Value = 1644  # M116 was 1832.5555
$ws.Cells.Item(116, 14).Value = -8964.1  # N116 was -9284

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 711950.7  # H132 was 792229.3
$ws.Cells.Item(132, 9).Value = 1618.1017  # I132 was 1689.9803
$ws.Cells.Item(132, 10).Value = 4902913  # J132 was 4457457.5
$ws.Cells.Item(132, 11).Value = 4854.3051  # K132 was 5069.9409
$ws.Cells.Item(132, 12).Value = 14708739  # L132 was 13372372.5
$ws.Cells.Item(132, 13).Value = -2324.3051  # M132 was -2539.9409
$ws.Cells.Item(132, 14).Value = -14713799  # N132 was -13377432.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 17861.4  # H135 was 20387.604
$ws.Cells.Item(135, 9).Value = 22504.848  # I135 was 25468.977
$ws.Cells.Item(135, 10).Value = 2604.3572  # J135 was 3026.25
$ws.Cells.Item(135, 11).Value = 202543.632  # K135 was 229220.793
$ws.Cells.Item(135, 12).Value = 23439.2148  # L135 was 27236.25
$ws.Cells.Item(135, 13).Value = -200008.632  # M135 was -226685.793
$ws.Cells.Item(135, 14).Value = -28509.2148  # N135 was -32306.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(136, 8).Value = 27639.75  # H136 was 27640
$ws.Cells.Item(136, 10).Value = 27639.75  # J136 was 27640
$ws.Cells.Item(136, 12).Value = 27639.75  # L136 was 27640
$ws.Cells.Item(136, 14).Value = -37839.75  # N136 was -37840

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2175375  # H137 was 2440561
$ws.Cells.Item(137, 9).Value = 3704873.2  # I137 was 4349131.5
$ws.Cells.Item(137, 10).Value = 1877.3684  # J137 was 1832.2222
$ws.Cells.Item(137, 11).Value = 11114619.6  # K137 was 13047394.5
$ws.Cells.Item(137, 12).Value = 5632.1052  # L137 was 5496.6666
$ws.Cells.Item(137, 13).Value = -11112069.6  # M137 was -13044844.5
$ws.Cells.Item(137, 14).Value = -10732.1052  # N137 was -10596.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1917387.9  # H138 was 2166426
$ws.Cells.Item(138, 9).Value = 1149.7046  # I138 was 1382.6666
$ws.Cells.Item(138, 10).Value = 3878189.5  # J138 was 4067439.8
$ws.Cells.Item(138, 11).Value = 3449.1138  # K138 was 4147.9998
$ws.Cells.Item(138, 12).Value = 11634568.5  # L138 was 12202319.4
$ws.Cells.Item(138, 13).Value = 1690.8862  # M138 was 992.0002000000004
$ws.Cells.Item(138, 14).Value = -11644848.5  # N138 was -12212599.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 93195  # H140 was 67778.17999999999
$ws.Cells.Item(140, 10).Value = 93195  # J140 was 67778.17999999999
$ws.Cells.Item(140, 12).Value = 93195  # L140 was 67778.17999999999
$ws.Cells.Item(140, 14).Value = -103555  # N140 was -78138.17999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 1528.8206  # H141 was 1304
$ws.Cells.Item(141, 9).Value = 1553.3055  # I141 was 1308.7046
$ws.Cells.Item(141, 11).Value = 4659.916499999999  # K141 was 3926.1138
$ws.Cells.Item(141, 13).Value = 520.0835000000006  # M141 was 1253.8862

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2587.3333  # H2 was 1542.4828
$ws.Cells.Item(2, 9).Value = 2422.2  # I2 was 902
$ws.Cells.Item(2, 10).Value = 2737.4546  # J2 was 2330.7693
$ws.Cells.Item(2, 11).Value = 2422.2  # K2 was 902
$ws.Cells.Item(2, 12).Value = 2737.4546  # L2 was 2330.7693
$ws.Cells.Item(2, 13).Value = -2309.2  # M2 was -789
$ws.Cells.Item(2, 14).Value = -2963.4546  # N2 was -2556.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 19115.75  # H36 was 8271
$ws.Cells.Item(36, 9).Value = 6585.2  # I36 was 8081.5
$ws.Cells.Item(36, 10).Value = 40000  # J36 was 9029
$ws.Cells.Item(36, 11).Value = 6585.2  # K36 was 8081.5
$ws.Cells.Item(36, 12).Value = 40000  # L36 was 9029
$ws.Cells.Item(36, 13).Value = -6239.2  # M36 was -7735.5
$ws.Cells.Item(36, 14).Value = -40692  # N36 was -9721

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 58942228  # H61 was 55667670
$ws.Cells.Item(61, 9).Value = 71500850  # I61 was 66734140
$ws.Cells.Item(61, 11).Value = 71500850  # K61 was 66734140
$ws.Cells.Item(61, 13).Value = -71500638  # M61 was -66733928

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 2718200  # H97 was 2718265.8
$ws.Cells.Item(97, 9).Value = 4465064.5  # I97 was 4167507.2
$ws.Cells.Item(97, 10).Value = 855.55554  # J97 was 937.5
$ws.Cells.Item(97, 11).Value = 4465064.5  # K97 was 4167507.2
$ws.Cells.Item(97, 12).Value = 855.55554  # L97 was 937.5
$ws.Cells.Item(97, 13).Value = -4464568.5  # M97 was -4167011.2
$ws.Cells.Item(97, 14).Value = -1847.55554  # N97 was -1929.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2587.3333  # H116 was 1542.4828
$ws.Cells.Item(116, 9).Value = 2422.2  # I116 was 902
$ws.Cells.Item(116, 10).Value = 2737.4546  # J116 was 2330.7693
$ws.Cells.Item(116, 11).Value = 2422.2  # K116 was 902
$ws.Cells.Item(116, 12).Value = 2737.4546  # L116 was 2330.7693
$ws.Cells.Item(116, 13).Value = -128.1999999999998  # M116 was 1392
$ws.Cells.Item(116, 14).Value = -7325.4546  # N116 was -6918.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 58942228  # H136 was 55667670
$ws.Cells.Item(136, 9).Value = 71500850  # I136 was 66734140
$ws.Cells.Item(136, 11).Value = 214502550  # K136 was 200202420
$ws.Cells.Item(136, 13).Value = -214500000  # M136 was -200199870

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2587.3333  # H3 was 1542.4828
$ws.Cells.Item(3, 9).Value = 2422.2  # I3 was 902
$ws.Cells.Item(3, 10).Value = 2737.4546  # J3 was 2330.7693
$ws.Cells.Item(3, 11).Value = 2422.2  # K3 was 902
$ws.Cells.Item(3, 12).Value = 2737.4546  # L3 was 2330.7693
$ws.Cells.Item(3, 13).Value = -2308.2  # M3 was -788
$ws.Cells.Item(3, 14).Value = -2965.4546  # N3 was -2558.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(33, 8).Value = 5001  # H33 was 4142.75
$ws.Cells.Item(33, 9).Value = 0  # I33 was 1713.3334
$ws.Cells.Item(33, 10).Value = 5001  # J33 was 5600.4
$ws.Cells.Item(33, 11).Value = 0  # K33 was 1713.3334
$ws.Cells.Item(33, 12).Value = 5001  # L33 was 5600.4
$ws.Cells.Item(33, 13).ClearContents()  # M33 was -1377.3334
$ws.Cells.Item(33, 14).Value = -5673  # N33 was -6272.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 433.08694  # H80 was 276.0909
$ws.Cells.Item(80, 10).Value = 523.06665  # J80 was 282.7857
$ws.Cells.Item(80, 12).Value = 523.06665  # L80 was 282.7857
$ws.Cells.Item(80, 14).Value = -2519.06665  # N80 was -2278.7857

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 433.08694  # H83 was 276.0909
$ws.Cells.Item(83, 10).Value = 523.06665  # J83 was 282.7857
$ws.Cells.Item(83, 12).Value = 2615.33325  # L83 was 1413.9285
$ws.Cells.Item(83, 14).Value = -12599.33325  # N83 was -11397.9285

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 1133.3334  # H19 was 303.83334
$ws.Cells.Item(19, 9).Value = 700  # I19 was 303.83334
$ws.Cells.Item(19, 10).Value = 2000  # J19 was 0
$ws.Cells.Item(19, 11).Value = 700  # K19 was 303.83334
$ws.Cells.Item(19, 12).Value = 2000  # L19 was 0
$ws.Cells.Item(19, 13).Value = -530  # M19 was -133.83334
$ws.Cells.Item(19, 14).Value = -2340  # N19 was empty

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(24, 8).Value = 1133.3334  # H24 was 303.83334
$ws.Cells.Item(24, 9).Value = 700  # I24 was 303.83334
$ws.Cells.Item(24, 10).Value = 2000  # J24 was 0
$ws.Cells.Item(24, 11).Value = 700  # K24 was 303.83334
$ws.Cells.Item(24, 12).Value = 2000  # L24 was 0
$ws.Cells.Item(24, 13).Value = -530  # M24 was -133.83334
$ws.Cells.Item(24, 14).Value = -2340  # N24 was empty

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 1433.3334  # H35 was 866.6667
$ws.Cells.Item(35, 9).Value = 900  # I35 was 866.6667
$ws.Cells.Item(35, 10).Value = 2500  # J35 was 0
$ws.Cells.Item(35, 11).Value = 900  # K35 was 866.6667
$ws.Cells.Item(35, 12).Value = 2500  # L35 was 0
$ws.Cells.Item(35, 13).Value = -606  # M35 was -572.6667
$ws.Cells.Item(35, 14).Value = -3088  # N35 was empty

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 25313.945  # H74 was 25314
$ws.Cells.Item(74, 10).Value = 25313.945  # J74 was 25314
$ws.Cells.Item(74, 12).Value = 25313.945  # L74 was 25314
$ws.Cells.Item(74, 14).Value = -27061.945  # N74 was -27062

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 25313.945  # H77 was 25314
$ws.Cells.Item(77, 10).Value = 25313.945  # J77 was 25314
$ws.Cells.Item(77, 12).Value = 75941.83499999999  # L77 was 75942
$ws.Cells.Item(77, 14).Value = -84677.83499999999  # N77 was -84678

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 36561.77  # H132 was 52671.258
$ws.Cells.Item(132, 9).Value = 23071.31  # I132 was 34853.367
$ws.Cells.Item(132, 10).Value = 91750  # J132 was 112064.22
$ws.Cells.Item(132, 11).Value = 69213.93000000001  # K132 was 104560.101
$ws.Cells.Item(132, 12).Value = 275250  # L132 was 336192.66
$ws.Cells.Item(132, 13).Value = -66683.93000000001  # M132 was -102030.101
$ws.Cells.Item(132, 14).Value = -280310  # N132 was -341252.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 16155.264  # H134 was 22602.234
$ws.Cells.Item(134, 9).Value = 1175.5  # I134 was 1434.9714
$ws.Cells.Item(134, 10).Value = 50200.184  # J134 was 68905.625
$ws.Cells.Item(134, 11).Value = 3526.5  # K134 was 4304.914199999999
$ws.Cells.Item(134, 12).Value = 150600.552  # L134 was 206716.875
$ws.Cells.Item(134, 13).Value = -991.5  # M134 was -1769.914199999999
$ws.Cells.Item(134, 14).Value = -155670.552  # N134 was -211786.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 28058  # H87 was 29283.03
$ws.Cells.Item(87, 9).Value = 20753.75  # I87 was 28653.5
$ws.Cells.Item(87, 10).Value = 28971.031  # J87 was 29366.967
$ws.Cells.Item(87, 11).Value = 62261.25  # K87 was 85960.5
$ws.Cells.Item(87, 12).Value = 86913.09299999999  # L87 was 88100.901
$ws.Cells.Item(87, 13).Value = -61013.25  # M87 was -84712.5
$ws.Cells.Item(87, 14).Value = -89409.09299999999  # N87 was -90596.901

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 28058  # H90 was 29283.03
$ws.Cells.Item(90, 9).Value = 20753.75  # I90 was 28653.5
$ws.Cells.Item(90, 10).Value = 28971.031  # J90 was 29366.967
$ws.Cells.Item(90, 11).Value = 186783.75  # K90 was 257881.5
$ws.Cells.Item(90, 12).Value = 260739.279  # L90 was 264302.703
$ws.Cells.Item(90, 13).Value = -180543.75  # M90 was -251641.5
$ws.Cells.Item(90, 14).Value = -273219.279  # N90 was -276782.703

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 1191.3334  # H92 was 1048.4
$ws.Cells.Item(92, 9).Value = 1037.25  # I92 was 1083
$ws.Cells.Item(92, 10).Value = 1499.5  # J92 was 996.5
$ws.Cells.Item(92, 11).Value = 3111.75  # K92 was 3249
$ws.Cells.Item(92, 12).Value = 4498.5  # L92 was 2989.5
$ws.Cells.Item(92, 13).Value = -1863.75  # M92 was -2001
$ws.Cells.Item(92, 14).Value = -6994.5  # N92 was -5485.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 3693.1667  # H97 was 2615.889
$ws.Cells.Item(97, 9).Value = 3693.1667  # I97 was 4395.8
$ws.Cells.Item(97, 10).Value = 0  # J97 was 391
$ws.Cells.Item(97, 11).Value = 11079.5001  # K97 was 13187.4
$ws.Cells.Item(97, 12).Value = 0  # L97 was 1173
$ws.Cells.Item(97, 13).Value = -10583.5001  # M97 was -12691.4
$ws.Cells.Item(97, 14).ClearContents()  # N97 was -2165

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 554.1842  # H113 was 881.63635
$ws.Cells.Item(113, 9).Value = 445.72726  # I113 was 500
$ws.Cells.Item(113, 10).Value = 1270  # J113 was 919.8
$ws.Cells.Item(113, 11).Value = 1337.18178  # K113 was 1500
$ws.Cells.Item(113, 12).Value = 3810  # L113 was 2759.4
$ws.Cells.Item(113, 13).Value = 832.8182200000001  # M113 was 670
$ws.Cells.Item(113, 14).Value = -8150  # N113 was -7099.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(126, 8).Value = 3188.889  # H126 was 2877.7778
$ws.Cells.Item(126, 9).Value = 600  # I126 was 1466.6666
$ws.Cells.Item(126, 10).Value = 3512.5  # J126 was 3583.3333
$ws.Cells.Item(126, 11).Value = 1800  # K126 was 4399.9998
$ws.Cells.Item(126, 12).Value = 10537.5  # L126 was 10749.9999
$ws.Cells.Item(126, 13).Value = 3140  # M126 was 540.0002000000004
$ws.Cells.Item(126, 14).Value = -20417.5  # N126 was -20629.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 3788.38  # H139 was 4215.154
$ws.Cells.Item(139, 9).Value = 2369.4211  # I139 was 2267.7896
$ws.Cells.Item(139, 10).Value = 4658.0645  # J139 was 5336.364
$ws.Cells.Item(139, 11).Value = 7108.263300000001  # K139 was 6803.3688
$ws.Cells.Item(139, 12).Value = 13974.1935  # L139 was 16009.092
$ws.Cells.Item(139, 13).Value = -1968.263300000001  # M139 was -1663.3688
$ws.Cells.Item(139, 14).Value = -24254.1935  # N139 was -26289.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 50000  # H104 was 0
$ws.Cells.Item(104, 10).Value = 50000  # J104 was 0
$ws.Cells.Item(104, 12).Value = 50000  # L104 was 0
$ws.Cells.Item(104, 14).Value = -56988  # N104 was empty

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2736.125  # H61 was 2990.2307
$ws.Cells.Item(61, 9).Value = 2679.3635  # I61 was 2919.2222
$ws.Cells.Item(61, 10).Value = 2861  # J61 was 3150
$ws.Cells.Item(61, 11).Value = 2679.3635  # K61 was 2919.2222
$ws.Cells.Item(61, 12).Value = 2861  # L61 was 3150
$ws.Cells.Item(61, 13).Value = -2477.3635  # M61 was -2717.2222
$ws.Cells.Item(61, 14).Value = -3265  # N61 was -3554

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 2736.125  # H113 was 2990.2307
$ws.Cells.Item(113, 9).Value = 2679.3635  # I113 was 2919.2222
$ws.Cells.Item(113, 10).Value = 2861  # J113 was 3150
$ws.Cells.Item(113, 11).Value = 2679.3635  # K113 was 2919.2222
$ws.Cells.Item(113, 12).Value = 2861  # L113 was 3150
$ws.Cells.Item(113, 13).Value = -509.3634999999999  # M113 was -749.2222000000002
$ws.Cells.Item(113, 14).Value = -7201  # N113 was -7490

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(137, 8).Value = 46205.5  # H137 was 57611.316
$ws.Cells.Item(137, 10).Value = 46205.5  # J137 was 57611.316
$ws.Cells.Item(137, 12).Value = 46205.5  # L137 was 57611.316
$ws.Cells.Item(137, 14).Value = -56405.5  # N137 was -67811.31599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 30000  # H139 was 0
$ws.Cells.Item(139, 10).Value = 30000  # J139 was 0
$ws.Cells.Item(139, 12).Value = 30000  # L139 was 0
$ws.Cells.Item(139, 14).Value = -40280  # N139 was empty

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(141, 8).Value = 75000  # H141 was 66798.664
$ws.Cells.Item(141, 10).Value = 75000  # J141 was 66798.664
$ws.Cells.Item(141, 12).Value = 75000  # L141 was 66798.664
$ws.Cells.Item(141, 14).Value = -85360  # N141 was -77158.664
